# Scheduled Moogle_Profits market-price refresh: update cached
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) on the
# affected leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 292.625
$ws.Range("J2").Value = 299
$ws.Range("L2").Value = 299
$ws.Range("N2").Value = -525
$ws.Range("H80").Value = 2792.7856
$ws.Range("J80").Value = 1812.8572
$ws.Range("L80").Value = 5438.571599999999
$ws.Range("N80").Value = -7434.571599999999
$ws.Range("H83").Value = 2792.7856
$ws.Range("J83").Value = 1812.8572
$ws.Range("L83").Value = 16315.7148
$ws.Range("N83").Value = -26299.7148
$ws.Range("H92").Value = 899.5
$ws.Range("I92").Value = 199
$ws.Range("K92").Value = 199
$ws.Range("M92").Value = 1049
$ws.Range("H137").Value = 2375.524
$ws.Range("I137").Value = 1954.125
$ws.Range("J137").Value = 3724
$ws.Range("K137").Value = 5862.375
$ws.Range("L137").Value = 11172
$ws.Range("M137").Value = -3312.375
$ws.Range("N137").Value = -16272

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2442.647
$ws.Range("I132").Value = 1979.5483
$ws.Range("K132").Value = 5938.644899999999
$ws.Range("M132").Value = -3408.644899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1004393.7
$ws.Range("I105").Value = 1823204.5
$ws.Range("K105").Value = 1823204.5
$ws.Range("M105").Value = -1821457.5
$ws.Range("H134").Value = 4596
$ws.Range("I134").Value = 3431.4827
$ws.Range("K134").Value = 10294.4481
$ws.Range("M134").Value = -7759.4481

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1396.3334
$ws.Range("I22").Value = 564.3333
$ws.Range("K22").Value = 564.3333
$ws.Range("M22").Value = -214.3333
$ws.Range("H31").Value = 7773.4473
$ws.Range("I31").Value = 3331.2307
$ws.Range("K31").Value = 3331.2307
$ws.Range("M31").Value = -3036.2307
$ws.Range("H34").Value = 7773.4473
$ws.Range("I34").Value = 3331.2307
$ws.Range("K34").Value = 3331.2307
$ws.Range("M34").Value = -3129.2307
$ws.Range("H52").Value = 135000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 135000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 135000
$ws.Range("M52").Value = ""
$ws.Range("N52").Value = -135588
$ws.Range("H58").Value = 5946.6895
$ws.Range("I58").Value = 5846.9287
$ws.Range("J58").Value = 6039.8
$ws.Range("K58").Value = 5846.9287
$ws.Range("L58").Value = 6039.8
$ws.Range("M58").Value = -5643.9287
$ws.Range("N58").Value = -6445.8
$ws.Range("H62").Value = 41668596
$ws.Range("I62").Value = 1375
$ws.Range("J62").Value = 50002044
$ws.Range("K62").Value = 1375
$ws.Range("L62").Value = 50002044
$ws.Range("M62").Value = -751
$ws.Range("N62").Value = -50003292
$ws.Range("H65").Value = 41668596
$ws.Range("I65").Value = 1375
$ws.Range("J65").Value = 50002044
$ws.Range("K65").Value = 6875
$ws.Range("L65").Value = 250010220
$ws.Range("M65").Value = -3755
$ws.Range("N65").Value = -250016460
$ws.Range("H68").Value = 60129
$ws.Range("J68").Value = 99990
$ws.Range("L68").Value = 99990
$ws.Range("N68").Value = -101488
$ws.Range("H71").Value = 60129
$ws.Range("J71").Value = 99990
$ws.Range("L71").Value = 299970
$ws.Range("N71").Value = -307458
$ws.Range("H105").Value = 1775.4286
$ws.Range("I105").Value = 1775.4286
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1775.4286
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -28.42859999999996
$ws.Range("N105").Value = ""
$ws.Range("H136").Value = 5946.6895
$ws.Range("I136").Value = 5846.9287
$ws.Range("J136").Value = 6039.8
$ws.Range("K136").Value = 17540.7861
$ws.Range("L136").Value = 18119.4
$ws.Range("M136").Value = -14990.7861
$ws.Range("N136").Value = -23219.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 10016.5
$ws.Range("J55").Value = 14749.75
$ws.Range("L55").Value = 44249.25
$ws.Range("N55").Value = -44603.25
$ws.Range("H60").Value = 319.2
$ws.Range("I60").Value = 204.75
$ws.Range("K60").Value = 614.25
$ws.Range("M60").Value = -363.25
$ws.Range("H97").Value = 1205.5294
$ws.Range("I97").Value = 866.8333
$ws.Range("K97").Value = 2600.4999
$ws.Range("M97").Value = -2104.4999
$ws.Range("H132").Value = 3072.9092
$ws.Range("I132").Value = 3200.3333
$ws.Range("K132").Value = 28802.9997
$ws.Range("M132").Value = -26272.9997
$ws.Range("H139").Value = 4414.9473
$ws.Range("I139").Value = 3227
$ws.Range("K139").Value = 9681
$ws.Range("M139").Value = -4541
$ws.Range("H140").Value = 1895.1428
$ws.Range("I140").Value = 1542.5714
$ws.Range("J140").Value = 2071.4285
$ws.Range("K140").Value = 4627.7142
$ws.Range("L140").Value = 6214.2855
$ws.Range("M140").Value = 552.2857999999997
$ws.Range("N140").Value = -16574.2855
$ws.Range("H141").Value = 4314.8
$ws.Range("I141").Value = 4314.8
$ws.Range("K141").Value = 12944.4
$ws.Range("M141").Value = -7764.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4945
$ws.Range("I102").Value = 3790.8333
$ws.Range("K102").Value = 3790.8333
$ws.Range("M102").Value = -2168.8333
$ws.Range("H132").Value = 5439.5435
$ws.Range("I132").Value = 4623.222
$ws.Range("K132").Value = 13869.666
$ws.Range("M132").Value = -11339.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5320.952
$ws.Range("I7").Value = 5263.3335
$ws.Range("J7").Value = 5666.6665
$ws.Range("K7").Value = 5263.3335
$ws.Range("L7").Value = 5666.6665
$ws.Range("M7").Value = -5151.3335
$ws.Range("N7").Value = -5890.6665
$ws.Range("H68").Value = 11264.087
$ws.Range("I68").Value = 9873.777
$ws.Range("J68").Value = 12157.857
$ws.Range("K68").Value = 9873.777
$ws.Range("L68").Value = 12157.857
$ws.Range("M68").Value = -9124.777
$ws.Range("N68").Value = -13655.857
$ws.Range("H69").Value = 50000
$ws.Range("I69").Value = 50000
$ws.Range("K69").Value = 50000
$ws.Range("M69").Value = -49189
$ws.Range("H71").Value = 11264.087
$ws.Range("I71").Value = 9873.777
$ws.Range("J71").Value = 12157.857
$ws.Range("K71").Value = 49368.885
$ws.Range("L71").Value = 60789.285
$ws.Range("M71").Value = -45624.885
$ws.Range("N71").Value = -68277.285
$ws.Range("H72").Value = 50000
$ws.Range("I72").Value = 50000
$ws.Range("K72").Value = 150000
$ws.Range("M72").Value = -145944
$ws.Range("H126").Value = 5320.952
$ws.Range("I126").Value = 5263.3335
$ws.Range("K126").Value = 15790.0005
$ws.Range("L126").Value = 16999.9995
$ws.Range("M126").Value = -13320.0005
$ws.Range("N126").Value = -21939.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3266.4482
$ws.Range("I132").Value = 2349.28
$ws.Range("K132").Value = 7047.84
$ws.Range("M132").Value = -4517.84

